# Appends new sensor-log rows to the PIR, Humidity, Temperature, and mmWave sheets,
# matching the 2026-01-28 16:50-16:51 batch recorded by the logger.
#
# Every column in the source log is stored as plain text. Excel's Range.Value
# setter auto-detects date-like strings ("2026-01-28") and percentage-like
# strings ("87.9%") and silently converts them to numbers, so those values are
# written with a leading apostrophe (the standard "force text" prefix; in a
# single-quoted PowerShell string a literal quote is written as two quotes,
# i.e. '''2026-01-28' is the single character ' followed by 2026-01-28).
# Plain strings (times, locations, statuses, "22.8C", etc.) round-trip as text
# without any extra escaping.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$startRow = 158
$rows = @(
    @('''2026-01-28', '16:50:11', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:13', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:18', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:23', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:28', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:33', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:38', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:43', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:48', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:53', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:50:58', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:51:03', '16:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('''2026-01-28', '16:51:08', '16:00', 'Bathroom', 'No Motion', 'Inactive')
)
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

$ws = $wb.Worksheets.Item("Humidity")
$startRow = 156
$rows = @(
    @('''2026-01-28', '16:50:10', '16:00', 'Bathroom', '''87.9%', 'Active'),
    @('''2026-01-28', '16:50:11', '16:00', 'Bathroom', '''86.9%', 'Active'),
    @('''2026-01-28', '16:50:18', '16:00', 'Bathroom', '''87.8%', 'Active'),
    @('''2026-01-28', '16:50:22', '16:00', 'Bathroom', '''87.8%', 'Active'),
    @('''2026-01-28', '16:50:30', '16:00', 'Bathroom', '''86.9%', 'Active'),
    @('''2026-01-28', '16:50:34', '16:00', 'Bathroom', '''87.8%', 'Active'),
    @('''2026-01-28', '16:50:38', '16:00', 'Bathroom', '''87.8%', 'Active'),
    @('''2026-01-28', '16:50:42', '16:00', 'Bathroom', '''86.9%', 'Active'),
    @('''2026-01-28', '16:50:50', '16:00', 'Bathroom', '''86.9%', 'Active'),
    @('''2026-01-28', '16:50:54', '16:00', 'Bathroom', '''87.9%', 'Active'),
    @('''2026-01-28', '16:50:58', '16:00', 'Bathroom', '''87.9%', 'Active'),
    @('''2026-01-28', '16:51:02', '16:00', 'Bathroom', '''86.9%', 'Active')
)
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

$ws = $wb.Worksheets.Item("Temperature")
$startRow = 156
$rows = @(
    @('''2026-01-28', '16:50:10', '16:00', 'Bathroom', '22.9C', 'Active'),
    @('''2026-01-28', '16:50:12', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:19', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:22', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:31', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:35', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:39', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:43', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:51', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:55', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:50:59', '16:00', 'Bathroom', '22.8C', 'Active'),
    @('''2026-01-28', '16:51:03', '16:00', 'Bathroom', '22.8C', 'Active')
)
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

$ws = $wb.Worksheets.Item("mmWave")
$startRow = 7
$rows = @(
    @('''2026-01-28', '16:50:57', '16:00', 'Living Room', 'NO_PRESENCE', 'Inactive'),
    @('''2026-01-28', '16:51:00', '16:00', 'Living Room', 'NO_PRESENCE', 'Inactive'),
    @('''2026-01-28', '16:51:04', '16:00', 'Living Room', 'NO_PRESENCE', 'Inactive'),
    @('''2026-01-28', '16:51:06', '16:00', 'Living Room', 'NO_PRESENCE', 'Inactive'),
    @('''2026-01-28', '16:51:09', '16:00', 'Living Room', 'NO_PRESENCE', 'Active')
)
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

